$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 19 (pushes existing rows 19-72 down to 20-73),
# then copy the now-shifted row 20 (the old row 19's data) into the
# new row 19 so it carries the same category/price info, and finally
# update its date (column D) to the new value from the diary update.
$ws.Rows("19:19").Insert()
$ws.Range("A20:R20").Copy()
$ws.Range("A19:R19").PasteSpecial()
$ws.Range("D19").Value = 44497

$excel.CutCopyMode = $false
